$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hydropower plant parameters")

# Insert a new row above row 17 ("Q_max_turb") to hold the new "no_turbines"
# parameter; Excel automatically shifts the rows below and re-points any
# formulas that referenced them.
$ws.Rows.Item(17).Insert()

$ws.Cells.Item(17, 1).Value = "no_turbines"
$ws.Cells.Item(17, 2).Value = "number of turbines (units)"
$ws.Cells.Item(17, 3).Value = 3
$ws.Cells.Item(17, 4).Value = 3

# The row-insert carries the formatting of the row above down onto the new
# row; reset it back to the plain/default style used by the other parameter
# rows.
$ws.Range("C17:D17").Style = "Normal"

# Reflect the author's final selection in the saved workbook.
$ws.Range("E17").Select() | Out-Null
